# gsc-export/Breadcrumbs.xlsx -- "updated main GSC export data"
#
# Appends the next day's row to the "Chart" sheet (A44:C44 = date,
# Invalid count, Valid count). The date must land as plain text (matching
# every other cell in column A), not get auto-converted to a date serial
# by Excel's type inference, and without pulling in a new number-format /
# style entry. The trick: write it as a text formula, then Copy +
# PasteSpecial(values) over itself to flatten it back down to a literal
# string cell with the sheet's default (unformatted) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = 44

$ws.Range("A$lastRow").Formula = '="2025-11-17"'
$ws.Range("A$lastRow").Copy()
$ws.Range("A$lastRow").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("B$lastRow").Value = 0
$ws.Range("C$lastRow").Value = 29
